# Apply the "fixed error in t-test sheet with v-lookup" edit described by the
# diff to Statistics/hypothesis_testing_amy.xlsx.
#
# Summary of the changes:
#  1. One_Sample_t_test sheet:
#       - Sample size (C5): 25 -> 20000
#       - Alpha (C9): 0.01 -> 0.0001 (1E-4)
#       - G18's VLOOKUP column index fixed: 7 -> 8 (the v-lookup bug fix
#         mentioned in the commit message)
#  2. two_sample_proportion_test sheet:
#       - Sample 1 data (D4 trials, D5 successes): 200/40 -> 103/11
#       - Sample 2 data (E4 trials, E5 successes): 200/22 -> 95/3
#       - Hypothesized difference (C7): 0.1 -> 0
#       - Type of test (C9): "2-sided" -> "1-sided Right"
#  3. Selection / active-sheet view state updated so that
#     one_sample_proportion_test becomes the active tab, with new
#     selected cells on each sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("One_Sample_t_test")
$ws2 = $wb.Worksheets.Item("one_sample_proportion_test")
$ws3 = $wb.Worksheets.Item("two_sample_proportion_test")

# --- One_Sample_t_test ---------------------------------------------------
$ws1.Range("C5").Value = 20000
$ws1.Range("C9").Value = 0.0001
$ws1.Range("G18").Formula = '=IF(VLOOKUP($C$10,$I$4:$P$6,8,0)<C9, "REJECT", "FAIL TO REJECT")'

# --- two_sample_proportion_test ------------------------------------------
$ws3.Range("D4").Value = 103
$ws3.Range("E4").Value = 95
$ws3.Range("D5").Value = 11
$ws3.Range("E5").Value = 3
$ws3.Range("C7").Value = 0
$ws3.Range("C9").Value = "1-sided Right"

# --- View / selection state ------------------------------------------------
# Update each sheet's selection; the last-selected sheet becomes the
# workbook's active tab, so select on one_sample_proportion_test last so
# that it ends up active (tabSelected / activeTab = 1), matching the target.
$ws1.Range("C5").Select()
$ws3.Range("D16").Select()
$ws2.Range("G27").Select()
